$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5

# Row 3
$ws.Range("C3").Value = "02:03"
$ws.Range("G3").Value = 3.1
$ws.Range("I3").Value = 2.25
$ws.Range("J3").Value = 3.75
$ws.Range("L3").Value = 2.87
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("U3").Value = 1.75
$ws.Range("V3").Value = 2
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 10
$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 11
$ws.Range("AJ3").Value = 9
$ws.Range("AP3").Value = 26
$ws.Range("AR3").Value = 81
$ws.Range("AU3").Value = 8

# Row 7
$ws.Range("G7").Value = 1.73
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.5
$ws.Range("K7").Value = 1.91
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("O7").Value = 1.57
$ws.Range("P7").Value = 2.25
$ws.Range("Q7").Value = 2.7
$ws.Range("R7").Value = 1.44
$ws.Range("S7").Value = 1.62
$ws.Range("T7").Value = 2.2
$ws.Range("U7").Value = 2.63
$ws.Range("V7").Value = 1.44
$ws.Range("Y7").Value = 10
$ws.Range("Z7").Value = 13
$ws.Range("AC7").Value = 6
$ws.Range("AH7").Value = 9.5
$ws.Range("AI7").Value = 23
$ws.Range("AO7").Value = 10
$ws.Range("AS7").Value = 351
$ws.Range("AT7").Value = 2.2
$ws.Range("AU7").Value = 11
$ws.Range("AV7").Value = 101
$ws.Range("AY7").Value = 51
$ws.Range("BA7").Value = 251

# Row 8
$ws.Range("N8").Value = 9
